$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (13) down onto the new row (14)
$ws.Range("A13:K13").Copy()
$ws.Range("A14:K14").PasteSpecial(-4122)

# Fill in the new prescale row values
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "L1_HTT120_SingleLLPJet40"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 1

# Update the active selection to match the saved view state
$ws.Range("C17").Select()
